$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7 (Generation 0-5): Fitness 7310 -> 7318
$ws.Range("C2:C7").Value = 7318

# Rows 8-252 (Generation 6-250): Fitness 7310 -> 7293
$ws.Range("C8:C252").Value = 7293
